# cryptos.xlsx price / 1h-volume refresh (GitHub Actions scheduled update).
# Coin #48/#49 (WEMIXTOKEN / Aave) also swapped rank positions this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.056.57'
$ws.Range('D3').Value = '1.922.51'
$ws.Range('E3').Value = '  +2.66%  '
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = "'326.28"
$ws.Range('E5').Value = '  +3.51%  '
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('D7').Value = "'0.5156"
$ws.Range('E7').Value = '  +1.67%  '
$ws.Range('D8').Value = "'0.3996"
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('D9').Value = "'0.08466"
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('D10').Value = "'42.83"
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('D12').Value = "'21.17"
$ws.Range('E12').Value = '  +3.86%  '
$ws.Range('D13').Value = "'6.323"
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').Value = '1.919.25'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = "'7.346"
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = "'1.001"
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').Value = "'94.32"
$ws.Range('E17').Value = '  +3.53%  '
$ws.Range('D18').Value = "'0.00001115"
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = "'0.06761"
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('D20').Value = "'17.98"
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').Value = "'6.052"
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('D23').Value = '30.070.02'
$ws.Range('E23').Value = '  +5.33%  '
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').Value = "'2.203"
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').Value = '2.139.80'
$ws.Range('E26').Value = '  +2.58%  '
$ws.Range('D27').Value = "'160.01"
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('D28').Value = "'20.98"
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('D29').Value = "'2.460"
$ws.Range('E29').Value = '  +4.32%  '
$ws.Range('D30').Value = "'128.98"
$ws.Range('E30').Value = '  +2.33%  '
$ws.Range('D31').Value = "'1.077"
$ws.Range('D32').Value = "'0.1056"
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('D33').Value = "'6.072"
$ws.Range('E33').Value = '  +4.89%  '
$ws.Range('D34').Value = "'3.653"
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('D35').Value = "'0.02497"
$ws.Range('E35').Value = '  +1.91%  '
$ws.Range('D36').Value = "'0.06604"
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Value = "'0.2216"
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('D38').Value = "'1.243"
$ws.Range('E38').Value = '  +4.51%  '
$ws.Range('D39').Value = "'8.999"
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').Value = "'5.193"
$ws.Range('E40').Value = '  +2.74%  '
$ws.Range('D41').Value = "'0.6524"
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').Value = "'1.241"
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D43').Value = "'11.43"
$ws.Range('E43').Value = '  +3.05%  '
$ws.Range('D44').Value = "'0.6130"
$ws.Range('E44').Value = '  +1.66%  '
$ws.Range('D45').Value = "'13.12"
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').Value = "'3.748"
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').Value = "'2.054"
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('D48').Value = "'1.242"
$ws.Range('E48').Value = '  +2.17%  '
$ws.Range('D49').Value = "'125.27"
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = "'79.30"
$ws.Range('E50').Value = '  +3.62%  '
$ws.Range('B51').Value = 'WEMIXTOKEN'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').Value = "'1.146"
$ws.Range('E51').Value = '  -2.59%  '
